$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.172.22'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.17%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.328.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.50%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '520.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.07%  '

$ws.Range('E7').Value = '  +0.23%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.538'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.43%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.354.61'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.49%  '

$ws.Range('E10').Value = '  +6.49%  '

$ws.Range('E11').Value = '  -0.83%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.23'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.92%  '

$ws.Range('E13').Value = '  +0.34%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.76%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.743.28'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.47%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '56.794.21'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.46%  '

$ws.Range('E17').Value = '  +2.25%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.339.16'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.76%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.49'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.14%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.01%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.35%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.30%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.06%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.81'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.32%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.162'
$ws.Range('D25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.994'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.26%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.42%  '

$ws.Range('E28').Value = '  +10.21%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.40%  '

$ws.Range('E30').Value = '  +5.22%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.72'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.78%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.65%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.31'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.10%  '

$ws.Range('E34').Value = '  +0.05%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.996'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.37%  '

$ws.Range('E36').Value = '  +1.59%  '

$ws.Range('E37').Value = '  +0.61%  '

$ws.Range('E38').Value = '  +3.73%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.55'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.67%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.90'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E41').Value = '  +0.34%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.59'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.41%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '137.02'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.92%  '

$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.20'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.78%  '

$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '278.57'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.32%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0935'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.57%  '

$ws.Range('E47').Value = '  +0.18%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.564'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.71%  '

$ws.Range('E49').Value = '  +4.44%  '

$ws.Range('E50').Value = '  +0.41%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.53'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.95%  '
